$excel = New-Object -ComObject Excel.Application
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy formatting (styles) from the last existing data row (43) down through the new rows (44-58)
$ws.Range("A43:V43").Copy()
$ws.Range("A44:V58").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 44
$ws.Range("A44").Value = 43
$ws.Range("B44").Value = "india"
$ws.Range("C44").Value = "isl"
$ws.Range("D44").Value = "2023-2024"
$ws.Range("E44").Value = 45263.64583333334
$ws.Range("F44").Value = "Goa"
$ws.Range("G44").Value = 1
$ws.Range("H44").Value = "Kerala Blasters"
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 1.86
$ws.Range("K44").Value = "29/11/2023 15:42"
$ws.Range("L44").Value = 1.83
$ws.Range("M44").Value = "03/12/2023 14:41"
$ws.Range("N44").Value = 3.67
$ws.Range("O44").Value = "29/11/2023 15:42"
$ws.Range("P44").Value = 3.89
$ws.Range("Q44").Value = "03/12/2023 15:29"
$ws.Range("R44").Value = 4.04
$ws.Range("S44").Value = "29/11/2023 15:42"
$ws.Range("T44").Value = 4.12
$ws.Range("U44").Value = "03/12/2023 14:40"
$ws.Range("V44").Value = "https://www.betexplorer.com/football/india/isl/fc-goa-kerala-blasters/6oDim7UR/"

# Row 45
$ws.Range("A45").Value = 44
$ws.Range("B45").Value = "india"
$ws.Range("C45").Value = "isl"
$ws.Range("D45").Value = "2023-2024"
$ws.Range("E45").Value = 45264.64583333334
$ws.Range("F45").Value = "East Bengal"
$ws.Range("G45").Value = 5
$ws.Range("H45").Value = "North East Utd"
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 1.78
$ws.Range("K45").Value = "27/11/2023 15:42"
$ws.Range("L45").Value = 2.16
$ws.Range("M45").Value = "04/12/2023 15:29"
$ws.Range("N45").Value = 3.81
$ws.Range("O45").Value = "27/11/2023 15:42"
$ws.Range("P45").Value = 3.6
$ws.Range("Q45").Value = "04/12/2023 15:29"
$ws.Range("R45").Value = 4.29
$ws.Range("S45").Value = "27/11/2023 15:42"
$ws.Range("T45").Value = 3.29
$ws.Range("U45").Value = "04/12/2023 15:29"
$ws.Range("V45").Value = "https://www.betexplorer.com/football/india/isl/east-bengal-north-east-united/8faJrRxq/"

# Row 46
$ws.Range("A46").Value = 45
$ws.Range("B46").Value = "india"
$ws.Range("C46").Value = "isl"
$ws.Range("D46").Value = "2023-2024"
$ws.Range("E46").Value = 45266.64583333334
$ws.Range("F46").Value = "Mohun Bagan"
$ws.Range("G46").Value = 2
$ws.Range("H46").Value = "Odisha FC"
$ws.Range("I46").Value = 2
$ws.Range("J46").Value = 1.76
$ws.Range("K46").Value = "05/12/2023 01:41"
$ws.Range("L46").Value = 1.71
$ws.Range("M46").Value = "06/12/2023 15:29"
$ws.Range("N46").Value = 3.91
$ws.Range("O46").Value = "05/12/2023 01:41"
$ws.Range("P46").Value = 4.18
$ws.Range("Q46").Value = "06/12/2023 15:29"
$ws.Range("R46").Value = 3.9
$ws.Range("S46").Value = "05/12/2023 01:41"
$ws.Range("T46").Value = 4.42
$ws.Range("U46").Value = "06/12/2023 15:25"
$ws.Range("V46").Value = "https://www.betexplorer.com/football/india/isl/mohun-bagan-odisha-fc/zHVFsohk/"

# Row 47
$ws.Range("A47").Value = 46
$ws.Range("B47").Value = "india"
$ws.Range("C47").Value = "isl"
$ws.Range("D47").Value = "2023-2024"
$ws.Range("E47").Value = 45267.64583333334
$ws.Range("F47").Value = "Jamshedpur"
$ws.Range("G47").Value = 2
$ws.Range("H47").Value = "Chennaiyin"
$ws.Range("I47").Value = 2
$ws.Range("J47").Value = 2.56
$ws.Range("K47").Value = "01/12/2023 15:42"
$ws.Range("L47").Value = 2.35
$ws.Range("M47").Value = "07/12/2023 15:20"
$ws.Range("N47").Value = 3.29
$ws.Range("O47").Value = "01/12/2023 15:42"
$ws.Range("P47").Value = 3.74
$ws.Range("Q47").Value = "07/12/2023 15:26"
$ws.Range("R47").Value = 2.64
$ws.Range("S47").Value = "01/12/2023 15:42"
$ws.Range("T47").Value = 2.84
$ws.Range("U47").Value = "07/12/2023 15:20"
$ws.Range("V47").Value = "https://www.betexplorer.com/football/india/isl/jamshedpur-chennaiyin-fc/E7UJt57e/"

# Row 48
$ws.Range("A48").Value = 47
$ws.Range("B48").Value = "india"
$ws.Range("C48").Value = "isl"
$ws.Range("D48").Value = "2023-2024"
$ws.Range("E48").Value = 45268.64583333334
$ws.Range("F48").Value = "Bengaluru FC"
$ws.Range("G48").Value = 0
$ws.Range("H48").Value = "Mumbai City"
$ws.Range("I48").Value = 4
$ws.Range("J48").Value = 3.41
$ws.Range("K48").Value = "01/12/2023 15:42"
$ws.Range("L48").Value = 3.37
$ws.Range("M48").Value = "08/12/2023 15:29"
$ws.Range("N48").Value = 3.5
$ws.Range("O48").Value = "01/12/2023 15:42"
$ws.Range("P48").Value = 3.7
$ws.Range("Q48").Value = "08/12/2023 15:29"
$ws.Range("R48").Value = 2.01
$ws.Range("S48").Value = "01/12/2023 15:42"
$ws.Range("T48").Value = 2.09
$ws.Range("U48").Value = "08/12/2023 15:29"
$ws.Range("V48").Value = "https://www.betexplorer.com/football/india/isl/bengaluru-fc-mumbai-city/8UTNuPM1/"

# Row 49
$ws.Range("A49").Value = 48
$ws.Range("B49").Value = "india"
$ws.Range("C49").Value = "isl"
$ws.Range("D49").Value = "2023-2024"
$ws.Range("E49").Value = 45269.64583333334
$ws.Range("F49").Value = "East Bengal"
$ws.Range("G49").Value = 0
$ws.Range("H49").Value = "Punjab"
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 1.81
$ws.Range("K49").Value = "08/12/2023 14:05"
$ws.Range("L49").Value = 1.89
$ws.Range("M49").Value = "09/12/2023 15:17"
$ws.Range("N49").Value = 3.5
$ws.Range("O49").Value = "08/12/2023 14:05"
$ws.Range("P49").Value = 3.61
$ws.Range("Q49").Value = "09/12/2023 15:17"
$ws.Range("R49").Value = 4.22
$ws.Range("S49").Value = "08/12/2023 14:05"
$ws.Range("T49").Value = 4.15
$ws.Range("U49").Value = "09/12/2023 15:17"
$ws.Range("V49").Value = "https://www.betexplorer.com/football/india/isl/east-bengal-minerva-punjab/MXXRvqx8/"

# Row 50
$ws.Range("A50").Value = 49
$ws.Range("B50").Value = "india"
$ws.Range("C50").Value = "isl"
$ws.Range("D50").Value = "2023-2024"
$ws.Range("E50").Value = 45270.64583333334
$ws.Range("F50").Value = "North East Utd"
$ws.Range("G50").Value = 1
$ws.Range("H50").Value = "Hyderabad"
$ws.Range("I50").Value = 1
$ws.Range("J50").Value = 3.4
$ws.Range("K50").Value = "04/12/2023 15:42"
$ws.Range("L50").Value = 2.23
$ws.Range("M50").Value = "10/12/2023 15:22"
$ws.Range("N50").Value = 3.46
$ws.Range("O50").Value = "04/12/2023 15:42"
$ws.Range("P50").Value = 3.39
$ws.Range("Q50").Value = "10/12/2023 15:22"
$ws.Range("R50").Value = 2.02
$ws.Range("S50").Value = "04/12/2023 15:42"
$ws.Range("T50").Value = 3.32
$ws.Range("U50").Value = "10/12/2023 15:22"
$ws.Range("V50").Value = "https://www.betexplorer.com/football/india/isl/north-east-united-hyderabad/Q5b8myt8/"

# Row 51
$ws.Range("A51").Value = 50
$ws.Range("B51").Value = "india"
$ws.Range("C51").Value = "isl"
$ws.Range("D51").Value = "2023-2024"
$ws.Range("E51").Value = 45272.64583333334
$ws.Range("F51").Value = "Goa"
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = "Mumbai City"
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 2.71
$ws.Range("K51").Value = "08/12/2023 15:43"
$ws.Range("L51").Value = 2.51
$ws.Range("M51").Value = "12/12/2023 15:28"
$ws.Range("N51").Value = 3.49
$ws.Range("O51").Value = "08/12/2023 15:43"
$ws.Range("P51").Value = 3.44
$ws.Range("Q51").Value = "12/12/2023 15:28"
$ws.Range("R51").Value = 2.39
$ws.Range("S51").Value = "08/12/2023 15:43"
$ws.Range("T51").Value = 2.82
$ws.Range("U51").Value = "12/12/2023 15:28"
$ws.Range("V51").Value = "https://www.betexplorer.com/football/india/isl/fc-goa-mumbai-city/4dcCneeE/"

# Row 52
$ws.Range("A52").Value = 51
$ws.Range("B52").Value = "india"
$ws.Range("C52").Value = "isl"
$ws.Range("D52").Value = "2023-2024"
$ws.Range("E52").Value = 45273.64583333334
$ws.Range("F52").Value = "Chennaiyin"
$ws.Range("G52").Value = 2
$ws.Range("H52").Value = "Bengaluru FC"
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 2.53
$ws.Range("K52").Value = "08/12/2023 15:43"
$ws.Range("L52").Value = 2.5
$ws.Range("M52").Value = "13/12/2023 15:21"
$ws.Range("N52").Value = 3.26
$ws.Range("O52").Value = "08/12/2023 15:43"
$ws.Range("P52").Value = 3.3
$ws.Range("Q52").Value = "13/12/2023 15:20"
$ws.Range("R52").Value = 2.69
$ws.Range("S52").Value = "08/12/2023 15:43"
$ws.Range("T52").Value = 2.93
$ws.Range("U52").Value = "13/12/2023 15:21"
$ws.Range("V52").Value = "https://www.betexplorer.com/football/india/isl/chennaiyin-fc-bengaluru-fc/W4QuVSPt/"

# Row 53
$ws.Range("A53").Value = 52
$ws.Range("B53").Value = "india"
$ws.Range("C53").Value = "isl"
$ws.Range("D53").Value = "2023-2024"
$ws.Range("E53").Value = 45274.64583333334
$ws.Range("F53").Value = "Punjab"
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = "Kerala Blasters"
$ws.Range("I53").Value = 1
$ws.Range("J53").Value = 3.44
$ws.Range("K53").Value = "13/12/2023 03:12"
$ws.Range("L53").Value = 3.15
$ws.Range("M53").Value = "14/12/2023 15:28"
$ws.Range("N53").Value = 3.33
$ws.Range("O53").Value = "13/12/2023 03:12"
$ws.Range("P53").Value = 3.29
$ws.Range("Q53").Value = "14/12/2023 15:28"
$ws.Range("R53").Value = 2.06
$ws.Range("S53").Value = "13/12/2023 03:12"
$ws.Range("T53").Value = 2.37
$ws.Range("U53").Value = "14/12/2023 15:28"
$ws.Range("V53").Value = "https://www.betexplorer.com/football/india/isl/minerva-punjab-kerala-blasters/AePqUnum/"

# Row 54
$ws.Range("A54").Value = 53
$ws.Range("B54").Value = "india"
$ws.Range("C54").Value = "isl"
$ws.Range("D54").Value = "2023-2024"
$ws.Range("E54").Value = 45275.64583333334
$ws.Range("F54").Value = "North East Utd"
$ws.Range("G54").Value = 1
$ws.Range("H54").Value = "Mohun Bagan"
$ws.Range("I54").Value = 3
$ws.Range("J54").Value = 3.48
$ws.Range("K54").Value = "13/12/2023 07:12"
$ws.Range("L54").Value = 3.88
$ws.Range("M54").Value = "15/12/2023 15:26"
$ws.Range("N54").Value = 3.44
$ws.Range("O54").Value = "13/12/2023 07:12"
$ws.Range("P54").Value = 3.79
$ws.Range("Q54").Value = "15/12/2023 15:29"
$ws.Range("R54").Value = 2.01
$ws.Range("S54").Value = "13/12/2023 07:12"
$ws.Range("T54").Value = 1.9
$ws.Range("U54").Value = "15/12/2023 15:29"
$ws.Range("V54").Value = "https://www.betexplorer.com/football/india/isl/north-east-united-mohun-bagan/hKElT6fg/"

# Row 55
$ws.Range("A55").Value = 54
$ws.Range("B55").Value = "india"
$ws.Range("C55").Value = "isl"
$ws.Range("D55").Value = "2023-2024"
$ws.Range("E55").Value = 45276.54166666666
$ws.Range("F55").Value = "Bengaluru FC"
$ws.Range("G55").Value = 1
$ws.Range("H55").Value = "Jamshedpur"
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 1.93
$ws.Range("K55").Value = "13/12/2023 15:42"
$ws.Range("L55").Value = 1.94
$ws.Range("M55").Value = "16/12/2023 12:50"
$ws.Range("N55").Value = 3.43
$ws.Range("O55").Value = "13/12/2023 15:42"
$ws.Range("P55").Value = 3.66
$ws.Range("Q55").Value = "16/12/2023 12:50"
$ws.Range("R55").Value = 3.75
$ws.Range("S55").Value = "13/12/2023 15:42"
$ws.Range("T55").Value = 3.87
$ws.Range("U55").Value = "16/12/2023 12:50"
$ws.Range("V55").Value = "https://www.betexplorer.com/football/india/isl/bengaluru-fc-jamshedpur/4fIhSQ9a/"

# Row 56
$ws.Range("A56").Value = 55
$ws.Range("B56").Value = "india"
$ws.Range("C56").Value = "isl"
$ws.Range("D56").Value = "2023-2024"
$ws.Range("E56").Value = 45276.64583333334
$ws.Range("F56").Value = "Mumbai City"
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = "East Bengal"
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 1.45
$ws.Range("K56").Value = "12/12/2023 15:42"
$ws.Range("L56").Value = 1.48
$ws.Range("M56").Value = "16/12/2023 15:26"
$ws.Range("N56").Value = 4.6
$ws.Range("O56").Value = "12/12/2023 15:42"
$ws.Range("P56").Value = 4.31
$ws.Range("Q56").Value = "16/12/2023 15:26"
$ws.Range("R56").Value = 5.63
$ws.Range("S56").Value = "12/12/2023 15:42"
$ws.Range("T56").Value = 7.03
$ws.Range("U56").Value = "16/12/2023 15:28"
$ws.Range("V56").Value = "https://www.betexplorer.com/football/india/isl/mumbai-city-east-bengal/bHHdRpP5/"

# Row 57
$ws.Range("A57").Value = 56
$ws.Range("B57").Value = "india"
$ws.Range("C57").Value = "isl"
$ws.Range("D57").Value = "2023-2024"
$ws.Range("E57").Value = 45277.64583333334
$ws.Range("F57").Value = "Odisha FC"
$ws.Range("G57").Value = 3
$ws.Range("H57").Value = "Hyderabad"
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 2.36
$ws.Range("K57").Value = "10/12/2023 15:42"
$ws.Range("L57").Value = 1.88
$ws.Range("M57").Value = "17/12/2023 15:26"
$ws.Range("N57").Value = 3.25
$ws.Range("O57").Value = "10/12/2023 15:42"
$ws.Range("P57").Value = 3.69
$ws.Range("Q57").Value = "17/12/2023 15:26"
$ws.Range("R57").Value = 2.91
$ws.Range("S57").Value = "10/12/2023 15:42"
$ws.Range("T57").Value = 4.11
$ws.Range("U57").Value = "17/12/2023 15:26"
$ws.Range("V57").Value = "https://www.betexplorer.com/football/india/isl/odisha-fc-hyderabad/x8G0Q4vC/"

# Row 58
$ws.Range("A58").Value = 57
$ws.Range("B58").Value = "india"
$ws.Range("C58").Value = "isl"
$ws.Range("D58").Value = "2023-2024"
$ws.Range("E58").Value = 45278.64583333334
$ws.Range("F58").Value = "Punjab"
$ws.Range("G58").Value = 1
$ws.Range("H58").Value = "Chennaiyin"
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 2.75
$ws.Range("K58").Value = "17/12/2023 15:12"
$ws.Range("L58").Value = 2.95
$ws.Range("M58").Value = "18/12/2023 15:27"
$ws.Range("N58").Value = 3.27
$ws.Range("O58").Value = "17/12/2023 15:12"
$ws.Range("P58").Value = 3.46
$ws.Range("Q58").Value = "18/12/2023 15:27"
$ws.Range("R58").Value = 2.47
$ws.Range("S58").Value = "17/12/2023 15:12"
$ws.Range("T58").Value = 2.4
$ws.Range("U58").Value = "18/12/2023 15:26"
$ws.Range("V58").Value = "https://www.betexplorer.com/football/india/isl/minerva-punjab-chennaiyin-fc/fqB5POgI/"
